# Workbook edit: rename MATCH_CARD_LINK -> MATCH_CODE and replace the full
# scorecard URL with the bare match code on both the "ODI Batting" and
# "ODI Bowling" sheets, then add a new "Player Info" sheet in front of them.
#
# NOTE: sheet references in this runtime resolve by current position, so all
# edits to the pre-existing sheets are made *before* a new sheet is inserted
# (which would otherwise shift everyone's index).

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 1. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code ---
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "4253"
    3 = "4254"
    4 = "4255"
    5 = "4256"
    6 = "4260"
    7 = "4285"
}
foreach ($row in $battingCodes.Keys) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
}

# --- 2. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCell = $bowlingSheet.Cells.Item(2, 2)
$bowlingCell.NumberFormat = "@"
$bowlingCell.Value = "4253"

# --- 3. Insert the new "Player Info" sheet in front of "ODI Batting" ---
# Re-fetch by current position (index 1) rather than relying on the
# previously-captured $battingSheet reference.
$infoSheet = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$infoSheet.Name = "Player Info"

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $infoHeaders.Length; $col++) {
    $cell = $infoSheet.Cells.Item(1, $col)
    $cell.Value = $infoHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$infoSheet.Cells.Item(2, 1).NumberFormat = "@"
$infoSheet.Cells.Item(2, 1).Value = "4840"
$infoSheet.Cells.Item(2, 2).Value = "John Dillon Campbell"
$infoSheet.Cells.Item(2, 3).Value = "Left Handed"
$infoSheet.Cells.Item(2, 4).Value = "Right Arm Off Break"
